$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "30.710.37"
Set-TextCell $ws "E2" "  +1.92%  "
Set-TextCell $ws "D3" "1.897.54"
Set-TextCell $ws "E3" "  +1.19%  "
Set-TextCell $ws "D4" "1.001"
Set-TextCell $ws "E4" "  +0.11%  "
Set-TextCell $ws "D5" "244.59"
Set-TextCell $ws "E5" "  +0.65%  "
Set-TextCell $ws "D6" "1.000"
Set-TextCell $ws "E6" "  +0.08%  "
Set-TextCell $ws "D7" "0.4932"
Set-TextCell $ws "E7" "  +0.59%  "
Set-TextCell $ws "D8" "0.2944"
Set-TextCell $ws "E8" "  +1.10%  "
Set-TextCell $ws "D9" "0.06777"
Set-TextCell $ws "E9" "  +2.81%  "
Set-TextCell $ws "D10" "1.894.83"
Set-TextCell $ws "E10" "  +1.02%  "
Set-TextCell $ws "D11" "17.28"
Set-TextCell $ws "E11" "  +5.32%  "
Set-TextCell $ws "D12" "0.07257"
Set-TextCell $ws "E12" "  +0.93%  "
Set-TextCell $ws "D13" "91.14"
Set-TextCell $ws "E13" "  +5.97%  "
Set-TextCell $ws "D14" "0.6787"
Set-TextCell $ws "E14" "  +1.85%  "
Set-TextCell $ws "D15" "5.043"
Set-TextCell $ws "E15" "  +2.63%  "
Set-TextCell $ws "D16" "30.705.38"
Set-TextCell $ws "E16" "  +2.05%  "
Set-TextCell $ws "D17" "0.000008032"
Set-TextCell $ws "E17" "  +3.09%  "
Set-TextCell $ws "D18" "1.000"
Set-TextCell $ws "E18" "  +0.08%  "
Set-TextCell $ws "D19" "13.14"
Set-TextCell $ws "E19" "  +2.92%  "
Set-TextCell $ws "D20" "2.140.69"
Set-TextCell $ws "E20" "  +0.76%  "
Set-TextCell $ws "D21" "1.001"
Set-TextCell $ws "E21" "  +0.35%  "
Set-TextCell $ws "D22" "4.815"
Set-TextCell $ws "E22" "  +0.97%  "
Set-TextCell $ws "D23" "196.47"
Set-TextCell $ws "E23" "  +36.70%  "
Set-TextCell $ws "D24" "6.111"
Set-TextCell $ws "E24" "  +4.96%  "
Set-TextCell $ws "D25" "9.402"
Set-TextCell $ws "E25" "  +2.51%  "
Set-TextCell $ws "D26" "156.41"
Set-TextCell $ws "E26" "  +2.32%  "
Set-TextCell $ws "D27" "19.26"
Set-TextCell $ws "E27" "  +13.58%  "
Set-TextCell $ws "D28" "1.907"
Set-TextCell $ws "E28" "  +0.92%  "
Set-TextCell $ws "D29" "1.399"
Set-TextCell $ws "E29" "  +0.24%  "
Set-TextCell $ws "D30" "4.319"
Set-TextCell $ws "E30" "  +2.70%  "
Set-TextCell $ws "E31" "  +3.71%  "
Set-TextCell $ws "D32" "4.010"
Set-TextCell $ws "E32" "  +0.49%  "
Set-TextCell $ws "D33" "0.05231"
Set-TextCell $ws "E33" "  +1.85%  "
Set-TextCell $ws "D34" "0.7438"
Set-TextCell $ws "E34" "  +3.92%  "
Set-TextCell $ws "D35" "1.111"
Set-TextCell $ws "E35" "  +0.27%  "
Set-TextCell $ws "D36" "2.756"
Set-TextCell $ws "E36" "  +3.43%  "
Set-TextCell $ws "D37" "0.01842"
Set-TextCell $ws "E37" "  -0.26%  "
Set-TextCell $ws "D38" "2.686"
Set-TextCell $ws "E38" "  +0.23%  "
Set-TextCell $ws "B39" "TrustWalletToken"
Set-TextCell $ws "C39" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws "D39" "0.9379"
Set-TextCell $ws "E39" "  +0.85%  "
Set-TextCell $ws "B40" "RenderToken"
Set-TextCell $ws "C40" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D40" "2.133"
Set-TextCell $ws "E40" "  -0.40%  "
Set-TextCell $ws "D41" "0.4421"
Set-TextCell $ws "E41" "  +4.70%  "
Set-TextCell $ws "D42" "105.40"
Set-TextCell $ws "E42" "  +1.79%  "
Set-TextCell $ws "B43" "PaxDollar"
Set-TextCell $ws "C43" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell $ws "D43" "1.001"
Set-TextCell $ws "E43" "  +0.21%  "
Set-TextCell $ws "B44" "FraxShare"
Set-TextCell $ws "C44" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D44" "5.757"
Set-TextCell $ws "E44" "  -0.45%  "
Set-TextCell $ws "D45" "7.582"
Set-TextCell $ws "E45" "  +2.86%  "
Set-TextCell $ws "D46" "0.1354"
Set-TextCell $ws "E46" "  +5.56%  "
Set-TextCell $ws "D47" "0.05864"
Set-TextCell $ws "E47" "  +3.15%  "
Set-TextCell $ws "D48" "8.774"
Set-TextCell $ws "E48" "  +6.85%  "
Set-TextCell $ws "B49" "NEARProtocol"
Set-TextCell $ws "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws "D49" "1.423"
Set-TextCell $ws "E49" "  +6.61%  "
Set-TextCell $ws "B50" "Decentraland"
Set-TextCell $ws "C50" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell $ws "D50" "0.3943"
Set-TextCell $ws "E50" "  +5.03%  "
Set-TextCell $ws "D51" "33.73"
Set-TextCell $ws "E51" "  +2.70%  "
